$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 343.2857
$ws.Range("I12").Value = 343.2857
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 343.2857
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -173.2857
$ws.Range("H17").Value = 2035.75
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2035.75
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6107.25
$ws.Range("N17").Value = -6443.25
$ws.Range("H33").Value = 2816.6667
$ws.Range("I33").Value = 2816.6667
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2816.6667
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -2587.6667
$ws.Range("H64").Value = 7361.524
$ws.Range("I64").Value = 5997.8184
$ws.Range("J64").Value = 8861.6
$ws.Range("K64").Value = 5997.8184
$ws.Range("L64").Value = 8861.6
$ws.Range("M64").Value = -5749.8184
$ws.Range("H67").Value = 7361.524
$ws.Range("I67").Value = 5997.8184
$ws.Range("J67").Value = 8861.6
$ws.Range("K67").Value = 5997.8184
$ws.Range("L67").Value = 8861.6
$ws.Range("M67").Value = -5139.8184
$ws.Range("H80").Value = 2739.1482
$ws.Range("I80").Value = 568.1818
$ws.Range("J80").Value = 4231.6875
$ws.Range("K80").Value = 1704.5454
$ws.Range("L80").Value = 12695.0625
$ws.Range("M80").Value = -706.5454
$ws.Range("N80").Value = -14691.0625
$ws.Range("H83").Value = 2739.1482
$ws.Range("I83").Value = 568.1818
$ws.Range("J83").Value = 4231.6875
$ws.Range("K83").Value = 5113.6362
$ws.Range("L83").Value = 38085.1875
$ws.Range("M83").Value = -121.6361999999999
$ws.Range("N83").Value = -48069.1875
$ws.Range("H93").Value = 30601
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 30601
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30601
$ws.Range("N93").Value = -35593
$ws.Range("H125").Value = 830.6
$ws.Range("I125").Value = 811.9091
$ws.Range("J125").Value = 882
$ws.Range("K125").Value = 7307.1819
$ws.Range("L125").Value = 7938
$ws.Range("M125").Value = -4847.1819
$ws.Range("N125").Value = -12858
$ws.Range("H126").Value = 74000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 74000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 74000
$ws.Range("N126").Value = -83880
$ws.Range("H135").Value = 17863416
$ws.Range("I135").Value = 25001232
$ws.Range("J135").Value = 18875
$ws.Range("K135").Value = 225011088
$ws.Range("L135").Value = 169875
$ws.Range("M135").Value = -225008553
$ws.Range("H137").Value = 2312.6667
$ws.Range("I137").Value = 2155.1482
$ws.Range("J137").Value = 3021.5
$ws.Range("K137").Value = 6465.444600000001
$ws.Range("L137").Value = 9064.5
$ws.Range("M137").Value = -3915.444600000001
$ws.Range("N137").Value = -14164.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 167
$ws.Range("I5").Value = 140
$ws.Range("J5").Value = 194
$ws.Range("K5").Value = 140
$ws.Range("L5").Value = 194
$ws.Range("M5").Value = -28
$ws.Range("H132").Value = 4030.0312
$ws.Range("I132").Value = 3866.5
$ws.Range("J132").Value = 4389.8
$ws.Range("K132").Value = 11599.5
$ws.Range("L132").Value = 13169.4
$ws.Range("M132").Value = -9069.5
$ws.Range("N132").Value = -18229.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 167
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 194
$ws.Range("K4").Value = 140
$ws.Range("L4").Value = 194
$ws.Range("M4").Value = -25
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H20").Value = 1953.2941
$ws.Range("I20").Value = 2043.0714
$ws.Range("J20").Value = 1534.3334
$ws.Range("K20").Value = 2043.0714
$ws.Range("L20").Value = 1534.3334
$ws.Range("M20").Value = -1796.0714
$ws.Range("N20").Value = -2028.3334
$ws.Range("H75").Value = 12208.728
$ws.Range("I75").Value = 7144.1113
$ws.Range("J75").Value = 34999.5
$ws.Range("K75").Value = 7144.1113
$ws.Range("L75").Value = 34999.5
$ws.Range("M75").Value = -6208.1113
$ws.Range("H78").Value = 12208.728
$ws.Range("I78").Value = 7144.1113
$ws.Range("J78").Value = 34999.5
$ws.Range("K78").Value = 21432.3339
$ws.Range("L78").Value = 104998.5
$ws.Range("M78").Value = -16752.3339
$ws.Range("H126").Value = 50000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 50000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880
$ws.Range("H130").Value = 50390
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 50390
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 50390
$ws.Range("N130").Value = -60430

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 199.5
$ws.Range("I22").Value = 199
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 199
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 151
$ws.Range("H92").Value = 36200.332
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 36200.332
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 36200.332
$ws.Range("N92").Value = -41192.332
$ws.Range("H109").Value = 42856.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 42856.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 42856.5
$ws.Range("N109").Value = -44936.5
$ws.Range("H132").Value = 2893.7896
$ws.Range("I132").Value = 2380.5386
$ws.Range("J132").Value = 4005.8333
$ws.Range("K132").Value = 7141.6158
$ws.Range("L132").Value = 12017.4999
$ws.Range("M132").Value = -4611.6158
$ws.Range("H134").Value = 2586.0476
$ws.Range("I134").Value = 1124.1666
$ws.Range("J134").Value = 6240.75
$ws.Range("K134").Value = 3372.4998
$ws.Range("L134").Value = 18722.25
$ws.Range("M134").Value = -837.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 163.92592
$ws.Range("I12").Value = 91.71429000000001
$ws.Range("J12").Value = 189.2
$ws.Range("K12").Value = 275.14287
$ws.Range("L12").Value = 567.5999999999999
$ws.Range("M12").Value = -102.14287
$ws.Range("N12").Value = -913.5999999999999
$ws.Range("H136").Value = 2299.2222
$ws.Range("I136").Value = 455
$ws.Range("J136").Value = 4143.4443
$ws.Range("K136").Value = 1365
$ws.Range("L136").Value = 12430.3329
$ws.Range("M136").Value = 3735
$ws.Range("N136").Value = -22630.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 416.7
$ws.Range("I13").Value = 336
$ws.Range("J13").Value = 470.5
$ws.Range("K13").Value = 336
$ws.Range("L13").Value = 470.5
$ws.Range("M13").Value = -197
$ws.Range("N13").Value = -748.5
$ws.Range("H70").Value = 71061.06
$ws.Range("I70").Value = 105600.27
$ws.Range("J70").Value = 7739.1665
$ws.Range("K70").Value = 105600.27
$ws.Range("L70").Value = 7739.1665
$ws.Range("M70").Value = -105330.27
$ws.Range("H73").Value = 71061.06
$ws.Range("I73").Value = 105600.27
$ws.Range("J73").Value = 7739.1665
$ws.Range("K73").Value = 105600.27
$ws.Range("L73").Value = 7739.1665
$ws.Range("M73").Value = -104664.27
$ws.Range("H80").Value = 72711.81
$ws.Range("I80").Value = 103999.55
$ws.Range("J80").Value = 3878.8
$ws.Range("K80").Value = 103999.55
$ws.Range("L80").Value = 3878.8
$ws.Range("M80").Value = -103001.55
$ws.Range("H83").Value = 72711.81
$ws.Range("I83").Value = 103999.55
$ws.Range("J83").Value = 3878.8
$ws.Range("K83").Value = 519997.75
$ws.Range("L83").Value = 19394
$ws.Range("M83").Value = -515005.75
$ws.Range("H97").Value = 436.64
$ws.Range("I97").Value = 354.9375
$ws.Range("J97").Value = 581.8889
$ws.Range("K97").Value = 354.9375
$ws.Range("L97").Value = 581.8889
$ws.Range("M97").Value = 141.0625
$ws.Range("N97").Value = -1573.8889
$ws.Range("H132").Value = 3206.3845
$ws.Range("I132").Value = 2878.2727
$ws.Range("J132").Value = 5011
$ws.Range("K132").Value = 8634.8181
$ws.Range("L132").Value = 15033
$ws.Range("M132").Value = -6104.8181
$ws.Range("N132").Value = -20093

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3371.2856
$ws.Range("I55").Value = 4400
$ws.Range("J55").Value = 799.5
$ws.Range("K55").Value = 4400
$ws.Range("L55").Value = 799.5
$ws.Range("M55").Value = -4227
$ws.Range("N55").Value = -1145.5
$ws.Range("H94").Value = 27500
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 27500
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 27500
$ws.Range("N94").Value = -28852
$ws.Range("H132").Value = 3349.805
$ws.Range("I132").Value = 2654.1516
$ws.Range("J132").Value = 6219.375
$ws.Range("K132").Value = 7962.4548
$ws.Range("L132").Value = 18658.125
$ws.Range("M132").Value = -5432.4548
$ws.Range("N132").Value = -23718.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 19655.25
$ws.Range("I37").Value = 14063
$ws.Range("J37").Value = 25247.5
$ws.Range("K37").Value = 14063
$ws.Range("L37").Value = 25247.5
$ws.Range("M37").Value = -13860
$ws.Range("H46").Value = 62633.332
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 62633.332
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 62633.332
$ws.Range("N46").Value = -63095.332
$ws.Range("H107").Value = 459.14285
$ws.Range("I107").Value = 492.5
$ws.Range("J107").Value = 445.8
$ws.Range("K107").Value = 1477.5
$ws.Range("L107").Value = 1337.4
$ws.Range("M107").Value = 442.5
$ws.Range("H126").Value = 1322.7858
$ws.Range("I126").Value = 1377.4584
$ws.Range("J126").Value = 994.75
$ws.Range("K126").Value = 4132.3752
$ws.Range("L126").Value = 2984.25
$ws.Range("M126").Value = -1662.3752
$ws.Range("N126").Value = -7924.25
$ws.Range("H132").Value = 1843.3088
$ws.Range("I132").Value = 1494.0656
$ws.Range("J132").Value = 4886.7144
$ws.Range("K132").Value = 4482.1968
$ws.Range("L132").Value = 14660.1432
$ws.Range("M132").Value = -1952.1968
$ws.Range("H134").Value = 62633.332
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 62633.332
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 187899.996
$ws.Range("N134").Value = -192969.996
